$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column I (Achievement Value | [NEW] | Photo taken? ...) ---
$ws.Columns("I:I").Insert()

# Approximate the inserted column's width/format off its left neighbour (column H),
# matching Excel's "insert column" behaviour of carrying formatting across.
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth

# --- New header cell ---
$ws.Cells.Item(1, 9).Value = "Twitter Username"

# --- Per-employee Twitter usernames (column I, rows 2-31) ---
$twitter = @{
  2  = "adamcogan"
  12 = "GregHarrisSSW"
  15 = "Jean_SSW"
  16 = "jernej_kavka"
  19 = "madkonst"
  20 = "liamelliott_au"
  21 = "matteightyate"
  22 = "MattGoldmanSSW"
  24 = "michaelsmedley"
  28 = "PennyWalker_SSW"
  30 = "ulyssesmac"
  31 = "William_DotNet"
}

foreach ($row in $twitter.Keys) {
  $ws.Cells.Item($row, 9).Value = $twitter[$row]
}

# --- Refresh the AutoFilter + named range so they cover the new column L ---
$ws.AutoFilterMode = $false
$ws.Range("A1:L31").AutoFilter()

foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$L`$31"
  }
}

# --- Restore view: scroll back to top-left and select I32 (matches author's final state) ---
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("I32").Select()
